# Weekly update: insert two new price observations for "Larga vida" tomatoes
# (Primera / Segunda) above the existing history, shifting all later rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(221).Insert()
$ws.Rows.Item(221).Insert()

function Set-Row($r, $vals) {
    $ws.Cells.Item($r,1).Value = $vals[0]
    $ws.Cells.Item($r,2).Value = $vals[1]
    $ws.Cells.Item($r,3).Value = $vals[2]
    $ws.Cells.Item($r,4).Value = $vals[3]
    $ws.Cells.Item($r,5).Value = $vals[4]
    $ws.Cells.Item($r,6).Value = $vals[5]
    $ws.Cells.Item($r,7).Value = $vals[6]
    $ws.Cells.Item($r,8).Value = $vals[7]
    $ws.Cells.Item($r,9).Value = $vals[8]
    $ws.Cells.Item($r,10).Value = $vals[9]
    $ws.Cells.Item($r,11).Value = $vals[10]
    $ws.Cells.Item($r,12).Value = $vals[11]
    $ws.Cells.Item($r,13).Value = $vals[12]
    $ws.Cells.Item($r,14).Value = $vals[13]
    $ws.Cells.Item($r,15).Value = $vals[14]
    $ws.Cells.Item($r,16).Value = $vals[15]
    $ws.Cells.Item($r,17).Value = $vals[16]
    $ws.Cells.Item($r,18).Value = $vals[17]
}

Set-Row 221 @("11","Vega Monumental Concepción","Bíobío",44474,8,100112020,"Tomate","Larga vida","Primera",600,19000,20000,19500,"`$/bandeja 18 kilos","Región de Arica y Parinacota",1083,18,"Hortaliza")
Set-Row 222 @("11","Vega Monumental Concepción","Bíobío",44474,8,100112020,"Tomate","Larga vida","Segunda",300,18000,18000,18000,"`$/bandeja 18 kilos","Región de Arica y Parinacota",1000,18,"Hortaliza")
